$wb = $excel.ActiveWorkbook

# ---- Sheet1 ----
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B43").Value = 6
$ws.Range("E43").Value = 8
$ws.Range("H43").Value = 10
$ws.Range("K43").Value = 4
$ws.Range("H43").Select()

# ---- Sheet2 ----
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Range("B43").Value = 7
$ws.Range("E43").Value = 16
$ws.Range("H43").Value = 9
$ws.Range("K43").Value = 8
$ws.Range("K43").Select()

# ---- Sheet3 ----
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("B43").Value = 8
$ws.Range("E43").Value = 15
$ws.Range("H43").Value = 20
$ws.Range("K43").Value = 12
$ws.Range("L43").Select()

# ---- Sheet4 ----
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Range("B43").Value = 7
$ws.Range("E43").Value = 7
$ws.Range("H43").Value = 13
$ws.Range("K43").Value = 14
$ws.Range("B44").Select()

# ---- Sheet5 ----
$ws = $wb.Worksheets.Item("Sheet5")
$ws.Range("B43").Value = 16
$ws.Range("E43").Value = 16
$ws.Range("H43").Value = 10
$ws.Range("K43").Value = 6
$ws.Range("L43").Select()

# ---- Sheet6 ----
$ws = $wb.Worksheets.Item("Sheet6")
$ws.Range("B43").Value = 30
$ws.Range("E43").Value = 22
$ws.Range("H43").Value = 16
$ws.Range("K43").Value = 12
$ws.Range("L43").Select()

# ---- Sheet7 ----
$ws = $wb.Worksheets.Item("Sheet7")
$ws.Range("K2").Value = 1
$ws.Range("H4").Value = 2
$ws.Range("B7").Value = 4
$ws.Range("E7").Value = 3
$ws.Range("K7").Value = 1
$ws.Range("H12").Value = 2
$ws.Range("K12").Value = 1
$ws.Range("B17").Value = 4
$ws.Range("E17").Value = 3
$ws.Range("K17").Value = 1
$ws.Range("H22").Value = 2
$ws.Range("K22").Value = 1
$ws.Range("E27").Value = 3
$ws.Range("K27").Value = 1
$ws.Range("H30").Value = 2
$ws.Range("K32").Value = 1
$ws.Range("B43").Value = 14
$ws.Range("E43").Value = 15
$ws.Range("H43").Value = 12
$ws.Range("K43").Value = 7
$ws.Range("I43").Select()

# ---- Sheet8 ----
$ws = $wb.Worksheets.Item("Sheet8")
$ws.Range("B2").Value = 5
$ws.Range("E2").Value = 5
$ws.Range("H2").Value = 5
$ws.Range("K2").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("K9").Value = 1
$ws.Range("B12").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("H12").Value = 5
$ws.Range("K12").Value = 5
$ws.Range("B13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("K22").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("B25").Value = 1
$ws.Range("K25").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("K28").Value = 1
$ws.Range("E29").Value = 1
$ws.Range("H30").Value = 1
$ws.Range("B31").Value = 1
$ws.Range("K31").Value = 1
$ws.Range("E32").Value = 1
$ws.Range("H33").Value = 1
$ws.Range("K34").Value = 1
$ws.Range("B43").Formula = "=11*2+9"
$ws.Range("E43").Value = 31
$ws.Range("H43").Value = 31
$ws.Range("K43").Value = 31
$ws.Range("H43").Select()

# ---- Sheet9 ----
$ws = $wb.Worksheets.Item("Sheet9")
$ws.Range("B2").Value = 5
$ws.Range("E2").Value = 4
$ws.Range("H2").Value = 3
$ws.Range("K2").Value = 2
$ws.Range("B7").Value = 5
$ws.Range("E7").Value = 4
$ws.Range("H7").Value = 3
$ws.Range("K7").Value = 2
$ws.Range("B12").Value = 4
$ws.Range("E12").Value = 3
$ws.Range("H12").Value = 2
$ws.Range("K12").Value = 1
$ws.Range("B17").Value = 4
$ws.Range("E17").Value = 3
$ws.Range("H17").Value = 2
$ws.Range("K17").Value = 1
$ws.Range("B22").Value = 3
$ws.Range("E22").Value = 2
$ws.Range("H22").Value = 1
$ws.Range("B27").Value = 3
$ws.Range("E27").Value = 2
$ws.Range("H27").Value = 1
$ws.Range("B32").Value = 2
$ws.Range("E32").Value = 1
$ws.Range("B43").Formula = "=22+14+10+3"
$ws.Range("E43").Formula = "=14+10+6+1"
$ws.Range("H43").Formula = "=10+6+2"
$ws.Range("K43").Value = 10
$ws.Range("O34").Select()

# ---- Sheet10 ----
$ws = $wb.Worksheets.Item("Sheet10")
$ws.Range("B2").Value = 5
$ws.Range("E2").Value = 5
$ws.Range("H2").Value = 5
$ws.Range("K2").Value = 5
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("B7").Value = 4
$ws.Range("E7").Value = 4
$ws.Range("H7").Value = 4
$ws.Range("K7").Value = 4
$ws.Range("B8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("B12").Value = 4
$ws.Range("E12").Value = 4
$ws.Range("H12").Value = 4
$ws.Range("K12").Value = 4
$ws.Range("B13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("B17").Value = 3
$ws.Range("E17").Value = 3
$ws.Range("H17").Value = 3
$ws.Range("K17").Value = 3
$ws.Range("B18").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("K18").Value = 1
$ws.Range("B22").Value = 3
$ws.Range("E22").Value = 3
$ws.Range("H22").Value = 3
$ws.Range("K22").Value = 3
$ws.Range("B23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("B27").Value = 2
$ws.Range("E27").Value = 2
$ws.Range("H27").Value = 2
$ws.Range("K27").Value = 2
$ws.Range("B28").Value = 1
$ws.Range("E28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("K28").Value = 1
$ws.Range("B32").Value = 2
$ws.Range("E32").Value = 2
$ws.Range("H32").Value = 2
$ws.Range("K32").Value = 2
$ws.Range("B33").Value = 1
$ws.Range("E33").Value = 1
$ws.Range("H33").Value = 1
$ws.Range("K33").Value = 1
$ws.Range("B43").Formula = "=11+14+10+6+7"
$ws.Range("E43").Value = 48
$ws.Range("H43").Value = 48
$ws.Range("K43").Value = 48
$ws.Range("O28").Select()
